$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($addr, $val) {
    $r = $ws.Range($addr)
    $r.NumberFormat = "@"
    $r.Value = $val
    $r.Style = "Normal"
}

Set-TextValue "D2" "29.910.41"
Set-TextValue "E2" "  +0.78%  "
Set-TextValue "D3" "1.635.73"
Set-TextValue "E3" "  +2.00%  "
Set-TextValue "E4" "  +0.48%  "
Set-TextValue "D5" "215.40"
Set-TextValue "E5" "  +1.42%  "
Set-TextValue "D6" "0.518"
Set-TextValue "E6" "  +0.39%  "
Set-TextValue "E7" "  +0.43%  "
Set-TextValue "D8" "28.90"
Set-TextValue "E8" "  +2.96%  "
Set-TextValue "D9" "0.258"
Set-TextValue "E9" "  +2.38%  "
Set-TextValue "D10" "0.0608"
Set-TextValue "E10" "  +0.90%  "
Set-TextValue "D11" "0.0915"
Set-TextValue "E11" "  +0.61%  "
Set-TextValue "D12" "1.872.45"
Set-TextValue "E12" "  +2.20%  "
Set-TextValue "D13" "1.640.46"
Set-TextValue "E13" "  +2.31%  "
Set-TextValue "D14" "0.568"
Set-TextValue "E14" "  +3.66%  "
Set-TextValue "D15" "9.29"
Set-TextValue "E15" "  +18.90%  "
Set-TextValue "D16" "3.88"
Set-TextValue "E16" "  +3.31%  "
Set-TextValue "D17" "29.971.67"
Set-TextValue "E17" "  +0.98%  "
Set-TextValue "D18" "64.25"
Set-TextValue "E18" "  +0.41%  "
Set-TextValue "D19" "244.22"
Set-TextValue "E19" "  +0.67%  "
Set-TextValue "D20" "0.0₃0704"
Set-TextValue "E20" "  +0.87%  "
Set-TextValue "E21" "  +0.31%  "
Set-TextValue "D22" "9.94"
Set-TextValue "E22" "  +5.35%  "
Set-TextValue "D23" "4.15"
Set-TextValue "E23" "  +3.04%  "
Set-TextValue "D24" "2.14"
Set-TextValue "E24" "  +1.77%  "
Set-TextValue "D25" "158.32"
Set-TextValue "E25" "  +1.90%  "
Set-TextValue "D26" "15.58"
Set-TextValue "E26" "  +0.54%  "
Set-TextValue "D27" "0.110"
Set-TextValue "E27" "  +1.69%  "
Set-TextValue "D28" "6.63"
Set-TextValue "E28" "  +2.88%  "
Set-TextValue "E29" "  +0.43%  "
Set-TextValue "D30" "0.0491"
Set-TextValue "E30" "  +1.79%  "
Set-TextValue "D31" "1.12"
Set-TextValue "E31" "  +5.37%  "
Set-TextValue "D32" "3.38"
Set-TextValue "E32" "  +4.37%  "
Set-TextValue "E33" "  -0.27%  "
Set-TextValue "D34" "1.428.10"
Set-TextValue "E34" "  -0.14%  "
Set-TextValue "D35" "1.66"
Set-TextValue "E35" "  +6.45%  "
Set-TextValue "D36" "1.04"
Set-TextValue "E36" "  +0.82%  "
Set-TextValue "E37" "  -2.35%  "
Set-TextValue "E38" "  +0.56%  "
Set-TextValue "D39" "0.0171"
Set-TextValue "E39" "  +1.39%  "
Set-TextValue "D40" "76.98"
Set-TextValue "E40" "  +15.90%  "
Set-TextValue "D41" "0.555"
Set-TextValue "E41" "  +1.57%  "
Set-TextValue "D42" "2.00"
Set-TextValue "E42" "  +1.70%  "
Set-TextValue "D43" "0.832"
Set-TextValue "E43" "  +1.74%  "
Set-TextValue "D44" "0.0494"
Set-TextValue "E44" "  -0.73%  "
Set-TextValue "B45" "WEMIXToken"
Set-TextValue "C45" "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
Set-TextValue "D45" "1.04"
Set-TextValue "E45" "  +7.03%  "
Set-TextValue "B46" "BitcoinSV"
Set-TextValue "C46" "https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv"
Set-TextValue "D46" "54.07"
Set-TextValue "E46" "  -7.44%  "
Set-TextValue "E47" "  +0.45%  "
Set-TextValue "D48" "5.37"
Set-TextValue "E48" "  +0.50%  "
Set-TextValue "D49" "1.778.94"
Set-TextValue "E49" "  +2.08%  "
Set-TextValue "D50" "0.0₆0112"
Set-TextValue "E50" "  +8.98%  "
Set-TextValue "D51" "89.56"
Set-TextValue "E51" "  +3.19%  "
